# "fixed export and fixing maps"
# Mtskheta Municipality area table: drop the census-data subtitle row and
# trim the historical year columns (1989 / 2002) down to just the latest
# (2014) figure, then restore the template's row heights / filler rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 held "(according to the population census data)" - remove it
# entirely so everything below shifts up one row.
$ws.Rows(2).Delete()

# Columns B:C held the 1989 and 2002 figures; delete them so the 2014
# column (old D) slides left into column B, taking its original
# right-edge border/number formatting with it.
$ws.Range("B:C").Delete()

# Re-apply the table's standard 20.1pt row height across the visible
# block, which also materializes the trailing blank rows 6-9.
for ($r = 1; $r -le 9; $r++) {
    $ws.Rows($r).RowHeight = 20.1
}
